# Update code from pycharm (#2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename process from "Update_Project_WellInfo" to "WellInfo_Update"
$ws.Range("B14").Value = "WellInfo_Update"

# Move active selection from B22 to E22
$ws.Range("E22").Select()
